# Rename the "network" / "location" header columns to "netid" / "name"
# (cleanup of name/location/loc column naming used by the dataframe readers).
#
# Set C2 ("network" -> "netid") before B2 ("location" -> "name") so that the
# new shared-string entries land in the same order as the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "netid"
$ws.Range("B2").Value = "name"

# Move the active selection to B3, matching the saved selection state.
$ws.Range("B3").Select() | Out-Null
